# The deck's single Design ("Integral", carried in ppt/theme/theme1.xml -
# the theme wired to the one-and-only Slide Master) is switched over to the
# built-in "Office Theme" colour palette (the palette that used to live in
# ppt/theme/theme2.xml, the Notes Master's theme).
#
# Office Theme's font scheme (Arial/Arial) and format scheme (fills, lines,
# effects, background fills) are byte-for-byte identical to Integral's, so
# the only real content change is the 10 accent/background colours (dk1 and
# lt1 - pure black/white - are shared by both palettes and stay the same).

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

function Set-ThemeColor {
    param($Scheme, [int] $Index, [string] $Hex)
    $r = [Convert]::ToInt32($Hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4,2), 16)
    # OLE COLORREF packs colours as 0x00BBGGRR
    $Scheme.Item($Index).RGB = ($b * 65536) + ($g * 256) + $r
}

# Office Theme colour scheme slots, in PowerPoint's fixed 1..12 order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
Set-ThemeColor $tcs 1  "000000"
Set-ThemeColor $tcs 2  "FFFFFF"
Set-ThemeColor $tcs 3  "44546A"
Set-ThemeColor $tcs 4  "E7E6E6"
Set-ThemeColor $tcs 5  "5B9BD5"
Set-ThemeColor $tcs 6  "ED7D31"
Set-ThemeColor $tcs 7  "A5A5A5"
Set-ThemeColor $tcs 8  "FFC000"
Set-ThemeColor $tcs 9  "4472C4"
Set-ThemeColor $tcs 10 "70AD47"
Set-ThemeColor $tcs 11 "0563C1"
Set-ThemeColor $tcs 12 "954F72"
